$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Harvest formatting for the future rows 19 & 20 from the current rows
# 21 & 22 ("Goals" / "Pass project ontime") before those rows are removed.
$ws.Range("A21:C21").Copy()
$ws.Range("A19:C19").PasteSpecial(-4122)
$ws.Range("A22:C22").Copy()
$ws.Range("A20:C20").PasteSpecial(-4122)

# --- Fix up the shading on what will become rows 12-13 (now unshaded,
# target wants the shaded "Medium" style used by row 14 onward) and row 5
# (currently shaded, target wants it unshaded like the rows above it).
$ws.Range("B4:C4").Copy()
$ws.Range("B5:C5").PasteSpecial(-4122)
$ws.Range("B14:C14").Copy()
$ws.Range("B12:C13").PasteSpecial(-4122)

# --- Remove the two now-redundant rows; this shifts the Q4/Yearly plan
# block up so the sheet ends at row 25 instead of row 27.
$ws.Rows("21:22").Delete()

# --- Rewrite the task list content ---
$ws.Range("A4").Value = 45567
$ws.Range("B4").Value = "Set up ssh-key using Azure Cloud Shell"

$ws.Range("B5").Value = "Create a virtual environment:"

$ws.Range("A6").ClearContents()
$ws.Range("B6").Value = "Run application in localhost"
$ws.Range("C6").Value = "Medium"

$ws.Range("B7").Value = "Create Github Actions"

$ws.Range("B8").Value = "Create Azure App Service"

$ws.Range("A9").Value = 45573
$ws.Range("B9").Value = "Create new agent pool"

$ws.Range("B10").Value = "Set up the Azure Pipeline"
$ws.Range("C10").Value = "Medium"

$ws.Range("A11").ClearContents()
$ws.Range("B11").Value = "Create new service connection"
$ws.Range("C11").ClearContents()

$ws.Range("A12").Value = 45581
$ws.Range("B12").Value = "Test and build pipeline"

$ws.Range("B13").Value = "Test app running on Azure"
$ws.Range("C13").Value = "Medium"

$ws.Range("A14").Value = "19/10/2024"
$ws.Range("B14").Value = "Capture the screenshots every step"

$ws.Range("B15").Value = "Update readme"
$ws.Range("C15").Value = "Easy"

$ws.Range("B16").Value = "Take screenshots demonstrating key steps"
$ws.Range("C16").Value = "Easy"

$ws.Range("B17").Value = "Create README file with instructions"
$ws.Range("C17").Value = "Easy"

$ws.Range("B18").Value = "Recording the video "

$ws.Range("B19").Value = "Goals"
$ws.Range("C19").ClearContents()

$ws.Range("B20").Value = "Pass project ontime"
$ws.Range("C20").ClearContents()

$ws.Range("B23").Value = "Pass this project and learning another course on Udacity"

# --- Restore the selection to match the saved workbook state ---
$ws.Range("B15").Select() | Out-Null
